$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A98").Value = "GRT-USD"
